# This script updates the "想去人数" (number of people who want to go) column (F)
# across the four worksheets of the Beijing comic-convention workbook, to reflect
# refreshed counts scraped at a later point in time (commit: "Update gh-pages to
# output generated at 456a3b4").

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibitions) ---
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F4").Value = 30
$ws.Range("F5").Value = 8178
$ws.Range("F8").Value = 92
$ws.Range("F9").Value = 7119
$ws.Range("F10").Value = 1135
$ws.Range("F11").Value = 549
$ws.Range("F12").Value = 492
$ws.Range("F14").Value = 710
$ws.Range("F18").Value = 231
$ws.Range("F21").Value = 78
$ws.Range("F22").Value = 11673
$ws.Range("F23").Value = 6
$ws.Range("F25").Value = 2293
$ws.Range("F27").Value = 3199
$ws.Range("F29").Value = 2706
$ws.Range("F31").Value = 25
$ws.Range("F32").Value = 286
$ws.Range("F35").Value = 1623
$ws.Range("F38").Value = 5833
$ws.Range("F39").Value = 81
$ws.Range("F40").Value = 1796
$ws.Range("F41").Value = 1246
$ws.Range("F42").Value = 848
$ws.Range("F43").Value = 162
$ws.Range("F44").Value = 188
$ws.Range("F45").Value = 1118
$ws.Range("F47").Value = 1528
$ws.Range("F48").Value = 101
$ws.Range("F49").Value = 1128

# --- Sheet "演出" (Performances) ---
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F8").Value = 251
$ws.Range("F15").Value = 6
$ws.Range("F16").Value = 108
$ws.Range("F20").Value = 67

# --- Sheet "本地生活" (Local life) ---
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 249
$ws.Range("F3").Value = 388

# --- Sheet "全部类型" (All types) ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 249
$ws.Range("F4").Value = 388
$ws.Range("F7").Value = 8178
$ws.Range("F10").Value = 92
$ws.Range("F11").Value = 7119
$ws.Range("F12").Value = 7119
$ws.Range("F13").Value = 1135
$ws.Range("F14").Value = 549
$ws.Range("F15").Value = 492
$ws.Range("F16").Value = 710
$ws.Range("F20").Value = 231
$ws.Range("F21").Value = 251
$ws.Range("F22").Value = 78
$ws.Range("F25").Value = 11673
$ws.Range("F27").Value = 6
$ws.Range("F29").Value = 2293
$ws.Range("F30").Value = 2293
$ws.Range("F31").Value = 3199
$ws.Range("F32").Value = 2706
$ws.Range("F33").Value = 25
$ws.Range("F34").Value = 286
$ws.Range("F38").Value = 1624
$ws.Range("F41").Value = 5833
$ws.Range("F42").Value = 67
$ws.Range("F43").Value = 1796
$ws.Range("F45").Value = 1246
$ws.Range("F46").Value = 848
$ws.Range("F47").Value = 188
$ws.Range("F48").Value = 1118
$ws.Range("F50").Value = 1528
$ws.Range("F51").Value = 1128
